# Apply the "adicionando alteracoes de media" edit:
#  - Replace the RANDBETWEEN formulas in column B (except rows 3 and 22) with
#    the static literal value 0 (formulas "frozen"/cleared to their resting value).
#  - Bump the row height of rows 1-22 from the LibreOffice default (13.8) to the
#    Excel default (15.75).
#  - Move the active selection from E27 to E12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose B-column RANDBETWEEN formula is replaced by a plain 0.
# (Rows 3 and 22 keep their original formula untouched.)
$rowsToClear = @(2,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21)

foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 2).Value = 0
}

# Rows 1-22 grow from 13.8pt to 15.75pt; rows 23-24 are left as-is.
$ws.Rows("1:22").RowHeight = 15.75

# Move the selection.
$ws.Range("E12").Select()
